$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# dominion_rps (column C) values updated from NULL to 1 for years 2046-2050 (rows 27-31)
$ws.Range("C27").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("C30").Value = 1
$ws.Range("C31").Value = 1

# Reflect the updated view/scroll position and selection left by the editor
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D25").Select()
